# "More drones and more batteris"
# - DJ Matrice 600 rows (9 & 10) now carry 6x the battery capacity
#   (modelled as a formula instead of a bare literal), which ripples
#   through the dependent Joules/Cruise/Power/PowerPerG formulas.
# - The "Battery Wh" label/value for the Penguin BE (rows 13-14) gets a
#   second, merged & centered column (C:D) to match the rest of the sheet.
# - Leave the cursor on I7, matching where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DJ Matrice 600 now has 6x the batteries -> battery mAh becomes a formula.
$ws.Range("C9").Formula = "=4500*6"
$ws.Range("C10").Formula = "=5700*6"

# Widen/merge the "Battery Wh" header and value cells over C:D and center them.
# (Set alignment first so both the header row and the value row only ever
# need one fresh style entry each, reused across C & D.)
$ws.Range("C13:D13").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C14:D14").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C13:D13").Merge()
$ws.Range("C14:D14").Merge()

# Restore the author's final selection.
$ws.Range("I7").Select()
